$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B ("ID Competição") holds 33 for every data row (rows 2-103); the
# competition id should be 233 instead - fix the dropped leading digit.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 103
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq 33) {
        $cell.Value = 233
    }
}
